$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header cell formatting (bold, border, alignment) from H1 into I1 and J1,
# then set the new header text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the plain data-cell formatting from column H into columns I and J
# for all data rows, then set the new numeric values.
$ws.Range("H2:H46").Copy($ws.Range("I2:J46"))

$data = @{
    2  = @(8, 8)
    3  = @(3, 4)
    4  = @(3, 4)
    5  = @(7, 7)
    6  = @(2, 2)
    7  = @(6, 7)
    8  = @(6, 7)
    9  = @(6, 6)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(8, 9)
    14 = @(3, 4)
    15 = @(1, 1)
    16 = @(7, 7)
    17 = @(10, 10)
    18 = @(4, 5)
    19 = @(6, 7)
    20 = @(5, 5)
    21 = @(8, 8)
    22 = @(7, 7)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(7, 7)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(9, 9)
    30 = @(11, 11)
    31 = @(6, 6)
    32 = @(1, 2)
    33 = @(5, 6)
    34 = @(4, 4)
    35 = @(7, 9)
    36 = @(7, 7)
    37 = @(5, 7)
    38 = @(6, 8)
    39 = @(7, 8)
    40 = @(7, 8)
    41 = @(5, 6)
    42 = @(2, 3)
    43 = @(1, 4)
    44 = @(4, 6)
    45 = @(3, 4)
    46 = @(2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
